$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "wYCsK239"
$ws.Range("B2").Value = 23071982
$ws.Range("C2").Value = "yhzsxcy35"
$ws.Range("D2").Value = 'k8$Qt2Z#'
$ws.Range("F2").Value = "czuFDRuU"
$ws.Range("G2").Value = "cafp"
